$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (22) of data to the table (name first, so the new shared
# string for the file name is registered before the new comment text)
$ws.Range("B22").Value = "1dmockanderrors10.csv"
$ws.Range("C22").Value = 400
$ws.Range("D22").Value = 50
$ws.Range("E22").Value = 3
$ws.Range("F22").Value = 0.3
$ws.Range("G22").Value = 60
$ws.Range("H22").Value = 1
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 20
$ws.Range("K22").Value = 20

# Update "Comments" column (L) for existing + new rows of the second table
# (the run counts were re-measured / re-labelled for the 2D interferograms
# with convolution work).
$ws.Range("L13").Value = "errorVSsamples done"
$ws.Range("L14").Value = "27 runs for 5%"
$ws.Range("L15").Value = "8 runs for 5%. New normalisation meathod. Ideal (coherent, noiseless) Interferogram peaks are now set to 1."
$ws.Range("L16").Value = "15 runs for 5%"
$ws.Range("L17").Value = "errorVSsamples done"
$ws.Range("L18").Value = "501 runs for 5%"
$ws.Range("L19").Value = "errorVSsamples done"
$ws.Range("L20").Value = "errorVSsamples done"
$ws.Range("L21").Value = "errorVSsamples done"
$ws.Range("L22").Value = "errorVSsamples done"

# Grow the table ("Table4") to include the newly added row
$lo = $ws.ListObjects.Item("Table4")
$lo.Resize($ws.Range("B12:L22"))

# Update the active selection to mirror the authored state
$ws.Range("L23").Select() | Out-Null
